# outputs/density_geounits/summary_stats_Lille.xlsx
# "create summary tables and figures"
#
# The "area_pop_sum" sheet currently has a "Density" column (C) duplicating
# the density figure already present in cell C2/C3, and a capitalized
# "Population" row label. Reshape it into a plain two-column index/value
# table: drop column C entirely, lowercase the "population" label, and add
# a new "density" row (row 4) carrying the value that used to live in C2/C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Drop column C ("Density" header + its two values) so the used range
# shrinks back down to columns A:B.
$ws.Range("C1:C3").Delete()

# "Population" -> "population"
$ws.Range("A3").Value = "population"

# New row 4: density label + value (same number that used to sit in C2/C3).
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 1859.380120108504
